$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Periodo Mora" values between the two data rows (E16/E17)
$ws.Range("E16").Value = "1902"
$ws.Range("E17").Value = "1903"

# Update "Valor Mora" (G16/G17) to the new amount
$ws.Range("G16").Value = 781242
$ws.Range("G17").Value = 781242
